$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39
$ws.Cells.Item($row, 1).Value = 38
$ws.Cells.Item($row, 2).Value = 51
$ws.Cells.Item($row, 3).Value = 9
$ws.Cells.Item($row, 4).Value = 14
$ws.Cells.Item($row, 5).Value = 17
$ws.Cells.Item($row, 6).Value = 74
$ws.Cells.Item($row, 7).Value = 91
